# gh-pages data refresh: "想去人数" (want-to-go headcount, column F) was
# re-scraped and the updated counts need to land on every sheet of the
# "杭州-漫展信息" workbook (the 全部类型 sheet mirrors rows from the other
# three, so several rows are touched twice under different sheet/row
# coordinates).
#
# $updates maps worksheet name -> { row number -> new column-F value } for
# every cell that changed between the previous scrape and this one.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2 = 1595; 3 = 3326; 4 = 30; 5 = 758; 6 = 2368; 7 = 508;
        8 = 427; 9 = 255; 10 = 153; 11 = 374; 12 = 1118; 13 = 470;
        14 = 230; 16 = 280; 17 = 4913; 19 = 1389; 20 = 3614; 21 = 163;
        22 = 213; 23 = 3921; 24 = 5284; 26 = 989; 27 = 579; 28 = 3384;
        29 = 395; 31 = 153; 33 = 905; 34 = 1224; 35 = 40; 36 = 62;
        37 = 1452; 38 = 148; 39 = 1434; 40 = 42; 41 = 930; 42 = 915;
        43 = 531; 44 = 64; 45 = 2458; 46 = 89; 47 = 188; 48 = 377;
        49 = 3765
    }
    "演出" = @{
        6 = 1033; 23 = 44; 26 = 47
    }
    "本地生活" = @{
        2 = 2604
    }
    "全部类型" = @{
        2 = 2604; 3 = 1595; 4 = 3326; 5 = 30; 6 = 758; 8 = 2368;
        9 = 508; 10 = 427; 11 = 255; 12 = 1033; 13 = 153; 14 = 374;
        15 = 1118; 16 = 470; 17 = 230; 19 = 280; 20 = 4913; 22 = 1389;
        23 = 3922; 24 = 5284; 26 = 989; 27 = 579; 28 = 3384; 29 = 395;
        31 = 153; 33 = 1224; 34 = 40; 35 = 62; 36 = 1452; 37 = 1434;
        38 = 930; 39 = 531; 41 = 64; 42 = 44; 43 = 2458; 45 = 89;
        46 = 188; 47 = 377; 49 = 3765
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowValues = $updates[$sheetName]
    foreach ($row in $rowValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowValues[$row]
    }
}
